# Commit: Fri, May 08, 2020 11:10:24 AM
#
# The presentation's design was changed from the "Integral" (Red Violet)
# theme to the stock "Office Theme" colour palette. In the underlying
# OOXML this shows up as the <a:clrScheme> colour slots of ppt/theme/theme1.xml
# (the theme used by the slide master / the whole deck) being swapped for the
# twelve standard Office theme colours.
#
# PowerPoint's automation model doesn't expose "swap the whole theme file" as
# a writable call (ApplyTheme only imports a *.thmx from disk, which isn't
# available in this environment) - the supported, persisted way to repaint a
# theme's colours from a script is via ThemeColorScheme.Colors(i).RGB, so we
# push each of the twelve standard "Office" RGB values into their slot, in
# the documented index order:
#   1 dk1, 2 lt1, 3 dk2, 4 lt2, 5 accent1, 6 accent2, 7 accent3,
#   8 accent4, 9 accent5, 10 accent6, 11 hlink, 12 folHlink

function Convert-RGBToLong([int]$r, [int]$g, [int]$b) {
    return $r + ($g * 256) + ($b * 65536)
}

$officeThemeColors = @(
    @(0x00, 0x00, 0x00),  # 1  dk1      #000000
    @(0xFF, 0xFF, 0xFF),  # 2  lt1      #FFFFFF
    @(0x44, 0x54, 0x6A),  # 3  dk2      #44546A
    @(0xE7, 0xE6, 0xE6),  # 4  lt2      #E7E6E6
    @(0x5B, 0x9B, 0xD5),  # 5  accent1  #5B9BD5
    @(0xED, 0x7D, 0x31),  # 6  accent2  #ED7D31
    @(0xA5, 0xA5, 0xA5),  # 7  accent3  #A5A5A5
    @(0xFF, 0xC0, 0x00),  # 8  accent4  #FFC000
    @(0x44, 0x72, 0xC4),  # 9  accent5  #4472C4
    @(0x70, 0xAD, 0x47),  # 10 accent6  #70AD47
    @(0x05, 0x63, 0xC1),  # 11 hlink    #0563C1
    @(0x95, 0x4F, 0x72)   # 12 folHlink #954F72
)

$p = $ppt.ActivePresentation

# The whole deck shares a single design/theme (one slide master), so grab it
# via the slide master's Theme and repaint every colour slot.
$themeColorScheme = $p.SlideMaster.Theme.ThemeColorScheme

for ($i = 1; $i -le $officeThemeColors.Count; $i++) {
    $rgb = $officeThemeColors[$i - 1]
    $themeColorScheme.Colors($i).RGB = Convert-RGBToLong $rgb[0] $rgb[1] $rgb[2]
}

Write-Host "Applied Office Theme colour scheme to the presentation theme."
